# Update the cryptocurrency price/volume table with the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.187.64"
$ws.Range("E2").Value = "  -0.38%  "
$ws.Range("D3").Value = "3.887.93"
$ws.Range("E3").Value = "  -1.03%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "483.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.96%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.622"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.744"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.80%  "
$ws.Range("E10").Value = "  +7.89%  "
$ws.Range("E11").Value = "  +0.64%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.20"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.64%  "
$ws.Range("E13").Value = "  -0.67%  "
$ws.Range("D14").Value = "4.503.31"
$ws.Range("E14").Value = "  -1.11%  "
$ws.Range("D15").Value = "3.908.72"
$ws.Range("E15").Value = "  -0.66%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.18"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.00%  "
$ws.Range("E17").Value = "  -0.57%  "
$ws.Range("E18").Value = "  +1.04%  "
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("D20").Value = "68.204.12"
$ws.Range("E20").Value = "  -0.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "430.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.58"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.92%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.80"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "89.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.77%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +16.70%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.66"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.83%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.32"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.65"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "711.92"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.50"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.87%  "
$ws.Range("E32").Value = "  +0.36%  "
$ws.Range("E33").Value = "  +2.76%  "
$ws.Range("D34").Value = "0.0₃0882"
$ws.Range("E34").Value = "  -0.87%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "61.51"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.20%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.05"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "40.90"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.146"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.40%  "
$ws.Range("B39").Value = "Dai"
$ws.Range("C39").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.998"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0499"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.394"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +14.80%  "
$ws.Range("B42").Value = "Fetch.AI"
$ws.Range("C42").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.07"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.37%  "
$ws.Range("E44").Value = "  -2.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.143"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.34"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.43%  "
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.36"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.52%  "
$ws.Range("E49").Value = "  -2.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "144.40"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.68%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₆0333"
$ws.Range("E51").Value = "  +23.23%  "
